# Insert a new data row above row 360. This pushes the existing rows
# 360..443 down to 361..444 (Excel's native Insert behavior also carries
# the row-360 formatting, e.g. the date style on column D, onto the new
# row), then we populate the newly-inserted row 360 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("360:360").Insert()

$ws.Range("A360").Value = 4
$ws.Range("B360").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C360").Value = 'Los Lagos'
$ws.Range("D360").Value = 45173
$ws.Range("E360").Value = 10
$ws.Range("F360").Value = 100112044
$ws.Range("G360").Value = 'Perejil'
$ws.Range("H360").Value = 'Sin especificar'
$ws.Range("I360").Value = 'Primera'
$ws.Range("J360").Value = 40
$ws.Range("K360").Value = 6000
$ws.Range("L360").Value = 6000
$ws.Range("M360").Value = 6000
$ws.Range("N360").Value = '$/docena de atados (3 kilos)'
$ws.Range("O360").Value = 'Región Metropolitana'
$ws.Range("P360").Value = 2000
$ws.Range("Q360").Value = 3
$ws.Range("R360").Value = 'Hortaliza'
